$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percentages, hour digits) are
# stored as text, matching the workbook author's original inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "328.91"
$ws.Range("E2").Value = "0.74%"
$ws.Range("G2").Value = "8"

$ws.Range("D3").Value = "43.89"
$ws.Range("E3").Value = "-1.66%"
$ws.Range("G3").Value = "8"

$ws.Range("D4").Value = "5.462"
$ws.Range("E4").Value = "-0.81%"
$ws.Range("G4").Value = "8"

$ws.Range("D5").Value = "0.07993"
$ws.Range("E5").Value = "-1.15%"
$ws.Range("G5").Value = "8"

$ws.Range("D6").Value = "2.002"
$ws.Range("E6").Value = "4.70%"
$ws.Range("G6").Value = "8"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.383"
$ws.Range("E7").Value = "2.22%"
$ws.Range("G7").Value = "8"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "2.588"
$ws.Range("E8").Value = "-3.89%"
$ws.Range("G8").Value = "8"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9499"
$ws.Range("E9").Value = "0.78%"
$ws.Range("G9").Value = "8"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1150"
$ws.Range("E10").Value = "-3.56%"
$ws.Range("G10").Value = "8"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1879"
$ws.Range("E11").Value = "0.68%"
$ws.Range("G11").Value = "8"

$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").Value = "11.70"
$ws.Range("E12").Value = "39.96%"
$ws.Range("G12").Value = "8"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "0.09939"
$ws.Range("E13").Value = "0.03%"
$ws.Range("G13").Value = "8"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.04831"
$ws.Range("E14").Value = "13.41%"
$ws.Range("G14").Value = "8"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.1064"
$ws.Range("E15").Value = "-0.29%"
$ws.Range("G15").Value = "8"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001264"
$ws.Range("E16").Value = "-1.08%"
$ws.Range("G16").Value = "8"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04072"
$ws.Range("E17").Value = "-2.74%"
$ws.Range("G17").Value = "8"

$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "0.005977"
$ws.Range("E18").Value = "-0.79%"
$ws.Range("G18").Value = "8"

$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "3.366"
$ws.Range("E19").Value = "-6.12%"
$ws.Range("G19").Value = "8"

$ws.Range("D20").Value = "0.3471"
$ws.Range("E20").Value = "-0.68%"
$ws.Range("G20").Value = "8"

$ws.Range("E21").Value = "3.77%"
$ws.Range("G21").Value = "8"

$ws.Range("E22").Value = "0.82%"
$ws.Range("G22").Value = "8"

$ws.Range("D23").Value = "0.001269"
$ws.Range("E23").Value = "2.62%"
$ws.Range("G23").Value = "8"

$ws.Range("D24").Value = "0.004334"
$ws.Range("E24").Value = "-4.49%"
$ws.Range("G24").Value = "8"

$ws.Range("D25").Value = "0.0001202"
$ws.Range("E25").Value = "1.80%"
$ws.Range("G25").Value = "8"

$ws.Range("D26").Value = "0.0003752"
$ws.Range("E26").Value = "-6.03%"
$ws.Range("G26").Value = "8"

$ws.Range("G27").Value = "8"

$ws.Range("G28").Value = "8"

$ws.Range("G29").Value = "8"

$ws.Range("G30").Value = "8"

$ws.Range("G31").Value = "8"

$ws.Range("G32").Value = "8"

$ws.Range("G33").Value = "8"

$ws.Range("G34").Value = "8"

$ws.Range("G35").Value = "8"

$ws.Range("G36").Value = "8"

$ws.Range("G37").Value = "8"

$ws.Range("D38").Value = "0.02566"
$ws.Range("E38").Value = "-2.66%"
$ws.Range("G38").Value = "8"

$ws.Range("D39").Value = "0.05637"
$ws.Range("E39").Value = "3.21%"
$ws.Range("G39").Value = "8"

$ws.Range("D40").Value = "0.007540"
$ws.Range("E40").Value = "-1.06%"
$ws.Range("G40").Value = "8"

$ws.Range("E41").Value = "-0.06%"
$ws.Range("G41").Value = "8"

$ws.Range("D42").Value = "0.007386"
$ws.Range("E42").Value = "3.05%"
$ws.Range("G42").Value = "8"

$ws.Range("D43").Value = "0.002019"
$ws.Range("E43").Value = "-0.30%"
$ws.Range("G43").Value = "8"

$ws.Range("D44").Value = "0.008614"
$ws.Range("E44").Value = "-2.41%"
$ws.Range("G44").Value = "8"

$ws.Range("D45").Value = "0.00007125"
$ws.Range("E45").Value = "0.28%"
$ws.Range("G45").Value = "8"

$ws.Range("D46").Value = "0.00000000752"
$ws.Range("E46").Value = "0.11%"
$ws.Range("G46").Value = "8"

$ws.Range("D47").Value = "0.003537"
$ws.Range("E47").Value = "55.69%"
$ws.Range("G47").Value = "8"

$ws.Range("D48").Value = "0.003791"
$ws.Range("E48").Value = "3.40%"
$ws.Range("G48").Value = "8"

$ws.Range("D49").Value = "0.00002104"
$ws.Range("E49").Value = "0.11%"
$ws.Range("G49").Value = "8"

$ws.Range("D50").Value = "0.0002004"
$ws.Range("E50").Value = "0.11%"
$ws.Range("G50").Value = "8"

$ws.Range("G51").Value = "8"
